$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we update to remain plain text so that
# values such as "1.000", "21.20", "0.000008522" keep their exact literal
# formatting instead of being auto-converted to numbers by Excel.
# (D43 is intentionally excluded - its value does not change.)
$ws.Range("D2:D42").NumberFormat = "@"
$ws.Range("D44:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.106.37"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.903.18"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "307.03"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "0.5229"
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("D8").Value = "0.3777"
$ws.Range("E8").Value = "  +3.08%  "
$ws.Range("D9").Value = "0.07219"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "21.20"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").Value = "0.8933"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "0.07667"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "1.895.31"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "94.44"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "5.235"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "0.000008522"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "14.52"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").Value = "0.9992"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "27.163.30"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").Value = "5.075"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "2.144.33"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").Value = "10.62"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").Value = "6.422"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "2.299"
$ws.Range("E25").Value = "  +9.62%  "
$ws.Range("D26").Value = "145.49"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").Value = "1.736"
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").Value = "18.10"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("D29").Value = "114.68"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "4.960"
$ws.Range("E30").Value = "  +4.40%  "
$ws.Range("D31").Value = "4.808"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("D32").Value = "0.09207"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "0.05065"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "1.241"
$ws.Range("E34").Value = "  +7.13%  "
$ws.Range("D35").Value = "0.7763"
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("D36").Value = "2.978"
$ws.Range("D37").Value = "3.300"
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").Value = "2.603"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "0.5662"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "0.01995"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").Value = "1.073"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "9.007"
$ws.Range("E42").Value = "  +5.15%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "119.04"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").Value = "0.1520"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("D46").Value = "0.4840"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D47").Value = "10.24"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("D48").Value = "0.9983"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "1.602"
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("D50").Value = "37.54"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").Value = "64.15"
$ws.Range("E51").Value = "  +1.56%  "
